# Updates the D (Price) and E (Volume(1h)) columns of the crypto
# tracking sheet with freshly scraped values, matching the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Every D/E cell in the sheet is stored as plain text (inline strings),
# including Price values that look like plain numbers (e.g. "1.00",
# "0.0220"). Assigning such a string straight to Range.Value lets Excel
# auto-detect it as a number and silently drop the formatting the feed
# relies on, so each such cell is switched to the Text number format
# ("@") first. Cells whose new text cannot be parsed as a number at all
# (e.g. "58.809.47", with two separators) do not need this - Excel
# already keeps those as text - so they are left untouched to avoid an
# unnecessary format change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D12",
    "D13",
    "D15",
    "D19",
    "D21",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
    )) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.809.47"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.496.71"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "536.95"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "136.97"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "2.519.84"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "2.942.66"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "22.99"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "58.733.00"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "2.509.65"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "11.11"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "323.88"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "65.77"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "7.54"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").Value = "6.70"
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "167.23"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "1.46"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "18.46"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "4.11"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").Value = "1.54"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("D39").Value = "36.69"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "284.67"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "133.05"
$ws.Range("E44").Value = "  +6.90%  "
$ws.Range("D45").Value = "0.994"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "10.88"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "0.0927"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "17.35"
$ws.Range("E51").Value = "  -3.05%  "
